# Insert a new data row before row 40 (pushes existing rows 40-119 down to 41-120),
# then populate the new row 40 with its own record.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows("40:40").Insert()

$ws.Range("A40").Value = 1
$ws.Range("B40").Value = "Agrícola del Norte S.A. de Arica"
$ws.Range("C40").Value = "Arica y Parinacota"
$ws.Range("D40").Value = 45246
$ws.Range("E40").Value = 15
$ws.Range("F40").Value = 100112012
$ws.Range("G40").Value = "Espinaca"
$ws.Range("H40").Value = "Sin especificar"
$ws.Range("I40").Value = "Tercera"
$ws.Range("J40").Value = 240
$ws.Range("K40").Value = 800
$ws.Range("L40").Value = 1000
$ws.Range("M40").Value = 925
$ws.Range("N40").Value = "`$/atado 2,5 a 3 kilos"
$ws.Range("O40").Value = "Región de Arica y Parinacota"
$ws.Range("P40").Value = 308
$ws.Range("Q40").Value = 3
$ws.Range("R40").Value = "Hortaliza"
